# Made distance scores less scale-dependent
# Update the raw counts on the "count" sheet; the "percent" sheet's
# formulas (=IF(count!Xn = 0, "", count!Xn/$Bn)) recalc automatically.

$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Item("count")
$percent = $wb.Worksheets.Item("percent")

# --- row 3 ---
$count.Range("D3").Value = 14653
$count.Range("E3").Value = 747

# --- row 4 ---
$count.Range("D4").Value = 732
$count.Range("E4").Value = 14668

# --- row 5 ---
$count.Range("F5").Value = 15694
$count.Range("J5").Value = 6

# --- row 6 ---
$count.Range("G6").Value = 15698
$count.Range("J6").Value = 2

# --- row 9 ---
$count.Range("G9").Value = 2
$count.Range("J9").Value = 15694

# --- row 10 ---
$count.Range("K10").Value = 10301
$count.Range("L10").Value = 5222
$count.Range("P10").Value = 177

# --- row 11 ---
$count.Range("K11").Value = 5125
$count.Range("L11").Value = 10381
$count.Range("P11").Value = 194

# --- row 12 ---
$count.Range("M12").Value = 11197
$count.Range("N12").Value = 4203

# --- row 13 ---
$count.Range("M13").Value = 4314
$count.Range("N13").Value = 11086

# --- row 15 ---
$count.Range("K15").Value = 188
$count.Range("L15").Value = 220
$count.Range("P15").Value = 15292

# --- row 16 ---
$count.Range("Q16").Value = 15343
$count.Range("R16").Value = 57

# --- row 17 ---
$count.Range("Q17").Value = 71
$count.Range("R17").Value = 15329

# --- row 20 ---
$count.Range("U20").Value = 15800
$count.Range("V20").Value = ""

# --- row 21 ---
$count.Range("U21").Value = 3
$count.Range("V21").Value = 15797

# --- row 28 ---
$count.Range("AC28").Value = 15782
$count.Range("AD28").Value = 18

# --- row 29 ---
$count.Range("AC29").Value = 23
$count.Range("AD29").Value = 15777

# Highlight the cells on "percent" whose mismatch rate moved enough to
# cross into the "flagged" (red-fill percent style) bucket.
$percent.Range("J5").Interior.Color = 255
$percent.Range("J6").Interior.Color = 255
$percent.Range("F9").Interior.Color = 255
$percent.Range("G9").Interior.Color = 255
$percent.Range("K15").Interior.Color = 255
$percent.Range("L15").Interior.Color = 255

# Restore selections / active sheet to match the saved view state.
$count.Activate()
$count.Range("A1:AP41").Select()

$percent.Activate()
$percent.Range("I12").Select()
